$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to remain plain text so Excel does not
# auto-convert numeric-looking strings (e.g. "119.88") into numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.486.46"
$ws.Range("D3").Value = "1.840.81"
$ws.Range("D5").Value = "333.99"
$ws.Range("D7").Value = "0.4621"
$ws.Range("D8").Value = "0.3849"
$ws.Range("D9").Value = "45.95"
$ws.Range("D10").Value = "0.07891"
$ws.Range("D11").Value = "0.9955"
$ws.Range("D12").Value = "21.41"
$ws.Range("D13").Value = "5.958"
$ws.Range("D14").Value = "1.844.05"
$ws.Range("D15").Value = "7.114"
$ws.Range("D16").Value = "1.007"
$ws.Range("D18").Value = "0.06684"
$ws.Range("D19").Value = "0.00001033"
$ws.Range("D21").Value = "1.006"
$ws.Range("D22").Value = "27.480.68"
$ws.Range("D23").Value = "5.380"
$ws.Range("D24").Value = "10.83"
$ws.Range("D25").Value = "2.312"
$ws.Range("D26").Value = "158.76"
$ws.Range("D27").Value = "2.053.99"
$ws.Range("D28").Value = "19.45"
$ws.Range("D29").Value = "2.106"
$ws.Range("D30").Value = "5.396"
$ws.Range("D31").Value = "119.88"
$ws.Range("D32").Value = "0.9728"
$ws.Range("D33").Value = "0.09384"
$ws.Range("D34").Value = "3.593"
$ws.Range("D35").Value = "5.289"
$ws.Range("D36").Value = "1.326"
$ws.Range("D37").Value = "0.06008"
$ws.Range("D38").Value = "0.02221"
$ws.Range("D40").Value = "1.178"
$ws.Range("D41").Value = "0.5886"
$ws.Range("D42").Value = "10.32"
$ws.Range("D43").Value = "0.1855"
$ws.Range("D44").Value = "1.237"
$ws.Range("D45").Value = "0.5569"
$ws.Range("D46").Value = "12.09"
$ws.Range("D47").Value = "1.904"
$ws.Range("D48").Value = "0.06682"
$ws.Range("D49").Value = "109.97"

# Restore the original (default) cell style on the Price column so no
# stray formatting differences are introduced.
$priceRange.Style = "Normal"

# Volume (E) column values are plain text already (they include
# surrounding spaces and a % sign) so they do not need special handling.
$ws.Range("E2").Value = "  -1.52%  "
$ws.Range("E3").Value = "  -2.10%  "
$ws.Range("E4").Value = "  -1.17%  "
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("E6").Value = "  -1.23%  "
$ws.Range("E7").Value = "  -0.67%  "
$ws.Range("E8").Value = "  -1.24%  "
$ws.Range("E9").Value = "  -2.03%  "
$ws.Range("E10").Value = "  -0.52%  "
$ws.Range("E11").Value = "  -0.98%  "
$ws.Range("E12").Value = "  -0.62%  "
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("E14").Value = "  -3.03%  "
$ws.Range("E15").Value = "  +0.78%  "
$ws.Range("E16").Value = "  -1.32%  "
$ws.Range("E17").Value = "  +1.74%  "
$ws.Range("E18").Value = "  -1.46%  "
$ws.Range("E19").Value = "  -0.59%  "
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("E21").Value = "  -1.20%  "
$ws.Range("E22").Value = "  -1.57%  "
$ws.Range("E23").Value = "  -1.35%  "
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("E25").Value = "  -1.65%  "
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("E27").Value = "  -3.07%  "
$ws.Range("E28").Value = "  -1.99%  "
$ws.Range("E29").Value = "  +2.67%  "
$ws.Range("E30").Value = "  -0.56%  "
$ws.Range("E31").Value = "  -0.79%  "
$ws.Range("E32").Value = "  +1.94%  "
$ws.Range("E33").Value = "  -0.78%  "
$ws.Range("E34").Value = "  -1.96%  "
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("E36").Value = "  -1.43%  "
$ws.Range("E37").Value = "  -1.61%  "
$ws.Range("E38").Value = "  -0.36%  "
$ws.Range("E39").Value = "  +2.27%  "
$ws.Range("E40").Value = "  -2.48%  "
$ws.Range("E41").Value = "  +0.34%  "
$ws.Range("E42").Value = "  +2.20%  "
$ws.Range("E43").Value = "  -1.23%  "
$ws.Range("E44").Value = "  -2.47%  "
$ws.Range("E45").Value = "  -0.71%  "
$ws.Range("E46").Value = "  +0.45%  "
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("E48").Value = "  -3.09%  "
$ws.Range("E49").Value = "  -2.98%  "
$ws.Range("E50").Value = "  -1.23%  "
$ws.Range("E51").Value = "  -1.39%  "
